$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("F2").Value = 40.66648600430862
$ws.Range("G2").Value = 39.88288365155705
$ws.Range("H2").Value = 41.45609404728989
$ws.Range("I2").Value = 0.001007105421425071
$ws.Range("J2").Value = 0.0007436697211170377
$ws.Range("K2").Value = 0.00150948627491243
$ws.Range("L2").Value = 0.0591696037836178
$ws.Range("M2").Value = 0.05813250605528451
$ws.Range("N2").Value = 0.06022991057654577

# Row 3
$ws.Range("F3").Value = 0.000009948548017113809
$ws.Range("G3").Value = 0.000000001969901810698318
$ws.Range("H3").Value = 0.00002952146945724846
$ws.Range("I3").Value = 0.00000836348374840006
$ws.Range("J3").Value = 0.000000001828839998814535
$ws.Range("K3").Value = 0.00002475954830356532
$ws.Range("L3").Value = 0.00001020524880153001
$ws.Range("M3").Value = 0.00000000206058932681671
$ws.Range("N3").Value = 0.00003027756615299297

# Row 4
$ws.Range("F4").Value = 40.66649595285663
$ws.Range("G4").Value = 39.88288365352695
$ws.Range("H4").Value = 41.45612356875935
$ws.Range("I4").Value = 0.001015468905173471
$ws.Range("J4").Value = 0.0007436715499570365
$ws.Range("K4").Value = 0.001534245823215995
$ws.Range("L4").Value = 0.05917980903241932
$ws.Range("M4").Value = 0.05813250811587384
$ws.Range("N4").Value = 0.06026018814269875
